$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("G2").Value = 0.053572
$ws.Range("H2").Value = 0.160716
$ws.Range("M2").Value = 6.322177333333333
$ws.Range("N2").Value = 18.966532
$ws.Range("O2").Value = 0.08271011762055308
$ws.Range("P2").Value = 0.08271011762055309
$ws.Range("Q2").Value = 0.3386916841013333
$ws.Range("R2").Value = 3.048225156912
$ws.Range("S2").Value = 0.08271011762055308
$ws.Range("T2").Value = 0.08271011762055309

# Row 3 (Target cluster: FAPs)
$ws.Range("G3").Value = 0.053572
$ws.Range("H3").Value = 0.160716
$ws.Range("O3").Value = 0.5401386314560596
$ws.Range("P3").Value = 0.5401386314560597
$ws.Range("Q3").Value = 2.21182689614
$ws.Range("R3").Value = 19.90644206526
$ws.Range("S3").Value = 0.5401386314560596
$ws.Range("T3").Value = 0.5401386314560597

# Row 4 (Target cluster: MuSCs)
$ws.Range("G4").Value = 0.053572
$ws.Range("H4").Value = 0.160716
$ws.Range("M4").Value = 27.73243066666667
$ws.Range("N4").Value = 83.197292
$ws.Range("O4").Value = 0.3628105447549136
$ws.Range("P4").Value = 0.3628105447549136
$ws.Range("Q4").Value = 1.485681775674667
$ws.Range("R4").Value = 13.371135981072
$ws.Range("S4").Value = 0.3628105447549136
$ws.Range("T4").Value = 0.3628105447549136

# Row 5 (Target cluster: Resolving-Mac)
$ws.Range("G5").Value = 0.053572
$ws.Range("H5").Value = 0.160716
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.096171666666667
$ws.Range("N5").Value = 3.288515
$ws.Range("O5").Value = 0.01434070616847367
$ws.Range("P5").Value = 0.01434070616847367
$ws.Range("Q5").Value = 0.05872410852666667
$ws.Range("R5").Value = 0.5285169767400001
$ws.Range("S5").Value = 0.01434070616847367
$ws.Range("T5").Value = 0.01434070616847367
